{"js": "// Replace the date line and each \"NNN\u00d7N=\" expression in the practice\n// table with its updated value. Every old value is unique within the\n// document, so a plain text search-and-replace (matching the whole\n// value, case-sensitive) unambiguously targets the correct run.\nconst replacements = [\n  [\"2024-03-11 Monday\", \"2024-03-12 Tuesday\"],\n  [\"173\u00d73=\", \"913\u00d76=\"],\n  [\"447\u00d74=\", \"309\u00d75=\"],\n  [\"691\u00d77=\", \"714\u00d77=\"],\n  [\"967\u00d78=\", \"425\u00d79=\"],\n  [\"747\u00d75=\", \"780\u00d77=\"],\n  [\"485\u00d79=\", \"786\u00d74=\"],\n  [\"993\u00d79=\", \"429\u00d73=\"],\n  [\"998\u00d77=\", \"902\u00d78=\"],\n  [\"577\u00d72=\", \"219\u00d79=\"],\n  [\"689\u00d72=\", \"209\u00d74=\"],\n  [\"949\u00d78=\", \"237\u00d72=\"],\n  [\"628\u00d73=\", \"834\u00d77=\"],\n  [\"987\u00d75=\", \"868\u00d75=\"],\n  [\"319\u00d76=\", \"250\u00d78=\"],\n  [\"369\u00d74=\", \"608\u00d75=\"],\n  [\"863\u00d79=\", \"653\u00d72=\"],\n  [\"785\u00d77=\", \"259\u00d76=\"],\n  [\"429\u00d72=\", \"196\u00d73=\"],\n  [\"554\u00d79=\", \"244\u00d73=\"],\n  [\"824\u00d73=\", \"686\u00d79=\"],\n  [\"968\u00d74=\", \"303\u00d72=\"],\n  [\"178\u00d79=\", \"112\u00d72=\"],\n  [\"210\u00d75=\", \"912\u00d77=\"],\n  [\"300\u00d78=\", \"388\u00d76=\"],\n  [\"377\u00d76=\", \"652\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00d7N=\" expression in the practice\n# table with its updated value. Every old value is unique within the\n# document, so a plain Find/Replace (case-sensitive, whole match) on\n# each pair unambiguously targets the correct run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-11 Monday\", \"2024-03-12 Tuesday\"),\n    @(\"173\u00d73=\", \"913\u00d76=\"),\n    @(\"447\u00d74=\", \"309\u00d75=\"),\n    @(\"691\u00d77=\", \"714\u00d77=\"),\n    @(\"967\u00d78=\", \"425\u00d79=\"),\n    @(\"747\u00d75=\", \"780\u00d77=\"),\n    @(\"485\u00d79=\", \"786\u00d74=\"),\n    @(\"993\u00d79=\", \"429\u00d73=\"),\n    @(\"998\u00d77=\", \"902\u00d78=\"),\n    @(\"577\u00d72=\", \"219\u00d79=\"),\n    @(\"689\u00d72=\", \"209\u00d74=\"),\n    @(\"949\u00d78=\", \"237\u00d72=\"),\n    @(\"628\u00d73=\", \"834\u00d77=\"),\n    @(\"987\u00d75=\", \"868\u00d75=\"),\n    @(\"319\u00d76=\", \"250\u00d78=\"),\n    @(\"369\u00d74=\", \"608\u00d75=\"),\n    @(\"863\u00d79=\", \"653\u00d72=\"),\n    @(\"785\u00d77=\", \"259\u00d76=\"),\n    @(\"429\u00d72=\", \"196\u00d73=\"),\n    @(\"554\u00d79=\", \"244\u00d73=\"),\n    @(\"824\u00d73=\", \"686\u00d79=\"),\n    @(\"968\u00d74=\", \"303\u00d72=\"),\n    @(\"178\u00d79=\", \"112\u00d72=\"),\n    @(\"210\u00d75=\", \"912\u00d77=\"),\n    @(\"300\u00d78=\", \"388\u00d76=\"),\n    @(\"377\u00d76=\", \"652\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n"}
